$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 2025-10-17 (serial 45947) day's data for both stations,
# mirroring the existing two-rows-per-day layout.
$ws.Range("A34").Value = 45947
$ws.Range("B34").Value = "四方坪站"
$ws.Range("C34").Value = 8361.11
$ws.Range("D34").Value = 6948.55
$ws.Range("E34").Value = 2890.07
$ws.Range("F34").Value = 373

$ws.Range("A35").Value = 45947
$ws.Range("B35").Value = "高岭站"
$ws.Range("C35").Value = 3470.37
$ws.Range("D35").Value = 2689.68
$ws.Range("E35").Value = 974.57
$ws.Range("F35").Value = 131

# Match the author's final selection/cursor position after entering the data.
$ws.Range("K31").Select()
